$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the "Expected Results" text of TC1 step 2.0 (D11) with an additional sentence
#    about also ordering the list by the date it arrived in the "liquidação" phase.
$ws.Range("D11").Value = "SYSTEM Exibe a lista de diárias (solicitações) aptas para pagamento ordenado pelo numero da diaria em ordem crescente. Exibe esta lista de diarias também ordenada pela data de chegada da solicitação na fase de liquidação (após registrar o empenho)."

# 2) Swap step 2.0 content between TC2 (rows 19-20) and TC3 (rows 27-28):
#    TC2's step becomes "Chefe Clica para realizar a liquidação." / "SYSTEM Apresenta a tela de Registrar Liquidações"
#    TC3's step becomes "Chefe Clica para atribuir/desatribuir o registro a si mesmo." / the related SYSTEM response.
$ws.Range("B20").Value = "Chefe Clica para realizar a liquidação."
$ws.Range("D20").Value = "SYSTEM Apresenta a tela de Registrar Liquidações"

$ws.Range("B28").Value = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
